# Update for ptdf: set `node_balance_sense` of all membering nodes of the
# ptdf_group to "none".
#
# This inserts a new "nodal_balance_sense" parameter column (with value
# "none" for every node row) into the obj_node-group_ptdf sheet, and makes
# that sheet the active/selected tab (it was previously
# rel_for_node_basic_structure that held the selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("obj_node-group_ptdf")

# Insert a new column before the existing "object_class"/ptdf-group columns
# (old column F), shifting the old F:J block to G:K.
$ws.Columns("F:F").Insert()

# New column takes on the same display width as the neighboring
# "para_name" column (E), matching Excel's normal insert-column behavior.
$ws.Columns("F:F").ColumnWidth = 23.83

# Header rows for the newly inserted column.
$ws.Range("F1").Value = "para_name"
$ws.Range("F2").Value = "nodal_balance_sense"

# Data rows: every node member of the ptdf_group gets "none".
$ws.Range("F4:F76").Value = "none"

# Make this sheet the active tab and restore the expected selection.
$ws.Activate()
$ws.Range("D3").Select()
